# Update "paises.xlsx" - countries & provincias Spain leaderboard refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 00:05"

# --- Update stats for rows whose country kept its row position ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1684092
$ws.Range("C4").Value = 17264
$ws.Range("D4").Value = 451451
$ws.Range("E4").Value = 1133398
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 560
$ws.Range("H4").Value = 99243

# Row 5: Brasil
$ws.Range("B5").Value = 360062
$ws.Range("C5").Value = 12664
$ws.Range("D5").Value = 142587
$ws.Range("E5").Value = 194872
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 590
$ws.Range("H5").Value = 22603

# Row 16: Canada
$ws.Range("B16").Value = 84699
$ws.Range("C16").Value = 1078
$ws.Range("D16").Value = 43985
$ws.Range("E16").Value = 34290
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 69
$ws.Range("H16").Value = 6424

# Row 25: Ecuador
$ws.Range("B25").Value = 36756
$ws.Range("C25").Value = 498
$ws.Range("D25").Value = 3560
$ws.Range("E25").Value = 30088
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 3108

# --- Colombia overtakes Ucrania in the ranking: swap rows 38/39 ---
# Row 38 becomes Colombia (with its newly updated stats)
$ws.Range("A38").Value = "Colombia"
$ws.Range("B38").Value = 21175
$ws.Range("C38").Value = 998
$ws.Range("D38").Value = 5016
$ws.Range("E38").Value = 15432
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 22
$ws.Range("H38").Value = 727

# Row 39 becomes Ucrania (its stats are unchanged from before)
$ws.Range("A39").Value = "Ucrania"
$ws.Range("B39").Value = 20986
$ws.Range("C39").Value = 406
$ws.Range("D39").Value = 7108
$ws.Range("E39").Value = 13261
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 12
$ws.Range("H39").Value = 617

# --- Guyana overtakes Bermudas & Islas Caimanes: shift rows 164-166 ---
# Row 164 becomes Guyana (with its newly updated stats)
$ws.Range("A164").Value = "Guyana"
$ws.Range("B164").Value = 135
$ws.Range("C164").Value = 8
$ws.Range("D164").Value = 62
$ws.Range("E164").Value = 63
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 10

# Row 165 becomes Bermudas (stats unchanged, shifted down one row)
$ws.Range("A165").Value = "Bermudas"
$ws.Range("B165").Value = 133
$ws.Range("C165").Value = 5
$ws.Range("D165").Value = 81
$ws.Range("E165").Value = 43
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 9

# Row 166 becomes Islas Caimanes (stats unchanged, shifted down one row)
$ws.Range("A166").Value = "Islas Caimanes"
$ws.Range("B166").Value = 129
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 61
$ws.Range("E166").Value = 67
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 1
